# fix lỗi trong report cơ sở. Thêm cột ghi chú trong báo cáo về chi tiêu
$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    # Assigning a DD-MM-YYYY-looking string directly makes Excel's COM
    # layer auto-convert it into a date serial (and stamp a date style
    # onto the cell). Forcing a Text number format first prevents the
    # conversion; clearing formats afterwards drops the now-unneeded
    # style so the cell is stored as a plain string, same as its
    # neighbours.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet "Đơn sale chính"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Đơn sale chính")

# Push the old "Tổng" row (row 3) down to row 5, freeing up rows 3-4 for
# two new service-order rows.
$ws1.Rows.Item(3).Insert()
$ws1.Rows.Item(4).Insert()

$ws1.Cells.Item(3,1).Value = "HD-LUXURY"
$ws1.Cells.Item(3,2).Value = 625
Set-TextValue $ws1.Cells.Item(3,3) "08-04-2024"
$ws1.Cells.Item(3,4).Value = "SÓC TRĂNG"
$ws1.Cells.Item(3,5).Value = "nguyễn thị mỹ chăm"
$ws1.Cells.Item(3,6).Value = "Cá nhân"
$ws1.Cells.Item(3,7).Value = "Cắt mí"
$ws1.Cells.Item(3,8).Value = 0
$ws1.Cells.Item(3,9).Value = "Lê Hoàng Thanh"
$ws1.Cells.Item(3,10).Value = 6000000
$ws1.Cells.Item(3,11).Value = 6000000
$ws1.Cells.Item(3,12).Value = 6000000
$ws1.Cells.Item(3,13).Value = 0.13
$ws1.Cells.Item(3,14).Value = 0

$ws1.Cells.Item(4,1).Value = "HD-LUXURY"
$ws1.Cells.Item(4,2).Value = 626
Set-TextValue $ws1.Cells.Item(4,3) "08-04-2024"
$ws1.Cells.Item(4,4).Value = "SÓC TRĂNG"
$ws1.Cells.Item(4,5).Value = "nguyễn thị mỹ trinh"
$ws1.Cells.Item(4,6).Value = "Cá nhân"
$ws1.Cells.Item(4,7).Value = "nhấn đồng tiền"
$ws1.Cells.Item(4,8).Value = 7000000
$ws1.Cells.Item(4,9).Value = 0
$ws1.Cells.Item(4,10).Value = 0
$ws1.Cells.Item(4,11).Value = 7000000
$ws1.Cells.Item(4,12).Value = 7000000
$ws1.Cells.Item(4,13).Value = 0
$ws1.Cells.Item(4,14).Value = 0

# Updated "Tổng" row, now row 5
$ws1.Cells.Item(5,2).Value = 3
$ws1.Cells.Item(5,8).Value = 16000000
$ws1.Cells.Item(5,10).Value = 6000000
$ws1.Cells.Item(5,11).Value = 22000000
$ws1.Cells.Item(5,12).Value = 22000000
$ws1.Cells.Item(5,13).Value = 0
$ws1.Cells.Item(5,14).Value = 900000

# ---------------------------------------------------------------------
# Sheet "Đơn 1 bác sĩ"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Đơn 1 bác sĩ")

$ws2.Rows.Item(3).Insert()
$ws2.Rows.Item(4).Insert()

$ws2.Cells.Item(3,1).Value = "HD-LUXURY"
$ws2.Cells.Item(3,2).Value = 625
Set-TextValue $ws2.Cells.Item(3,3) "08-04-2024"
$ws2.Cells.Item(3,4).Value = "SÓC TRĂNG"
$ws2.Cells.Item(3,5).Value = "nguyễn thị mỹ chăm"
$ws2.Cells.Item(3,6).Value = "Cá nhân"
$ws2.Cells.Item(3,7).Value = "Cắt mí"
$ws2.Cells.Item(3,8).Value = 0
$ws2.Cells.Item(3,9).Value = "Lê Hoàng Thanh"
$ws2.Cells.Item(3,10).Value = 6000000
$ws2.Cells.Item(3,11).Value = 6000000
$ws2.Cells.Item(3,12).Value = 6000000
$ws2.Cells.Item(3,13).Value = 0.1
$ws2.Cells.Item(3,14).Value = 600000

$ws2.Cells.Item(4,1).Value = "HD-LUXURY"
$ws2.Cells.Item(4,2).Value = 627
Set-TextValue $ws2.Cells.Item(4,3) "08-04-2024"
$ws2.Cells.Item(4,4).Value = "SÓC TRĂNG"
Set-TextValue $ws2.Cells.Item(4,5) "tạ duy hoàng "
$ws2.Cells.Item(4,6).Value = "Cá nhân"
$ws2.Cells.Item(4,7).Value = "Cắt mí"
$ws2.Cells.Item(4,8).Value = 6000000
$ws2.Cells.Item(4,9).Value = 0
$ws2.Cells.Item(4,10).Value = 0
$ws2.Cells.Item(4,11).Value = 6000000
$ws2.Cells.Item(4,12).Value = 6000000
$ws2.Cells.Item(4,13).Value = 0.1
$ws2.Cells.Item(4,14).Value = 600000

# Updated "Tổng" row, now row 5
$ws2.Cells.Item(5,2).Value = 3
$ws2.Cells.Item(5,8).Value = 14000000
$ws2.Cells.Item(5,10).Value = 6000000
$ws2.Cells.Item(5,11).Value = 20000000
$ws2.Cells.Item(5,12).Value = 18000000
$ws2.Cells.Item(5,13).Value = 0
$ws2.Cells.Item(5,14).Value = 1800000

# ---------------------------------------------------------------------
# Sheet "Đơn thu nợ"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Đơn thu nợ")

# Push the old "Tổng" row (row 4) down to row 5, freeing up row 4 for a
# new debt-collection row.
$ws3.Rows.Item(4).Insert()

$ws3.Cells.Item(4,1).Value = "TN"
$ws3.Cells.Item(4,2).Value = 181
$ws3.Cells.Item(4,3).Value = 1500000
Set-TextValue $ws3.Cells.Item(4,4) "08-09-2024"
$ws3.Cells.Item(4,5).Value = "CẦN THƠ"
$ws3.Cells.Item(4,6).Value = "HD-LUXURY-538"
$ws3.Cells.Item(4,7).Value = "Nâng mũi"
$ws3.Cells.Item(4,8).Value = "Ngô Xuân Nhi"
$ws3.Cells.Item(4,9).Value = "Cá nhân"
$ws3.Cells.Item(4,10).Value = "Lâm Hoàng Phú"
$ws3.Cells.Item(4,11).Value = 10000000
$ws3.Cells.Item(4,12).Value = "Đỗ Thị Huyền Trân"
$ws3.Cells.Item(4,13).Value = 8000000
$ws3.Cells.Item(4,14).Value = 18000000
$ws3.Cells.Item(4,15).Value = 11000000
$ws3.Cells.Item(4,16).Value = "Lâm Thị Mỹ Hằng"
$ws3.Cells.Item(4,17).Value = 0
$ws3.Cells.Item(4,18).Value = 0
$ws3.Cells.Item(4,19).Value = 0
$ws3.Cells.Item(4,20).Value = 0
$ws3.Cells.Item(4,21).Value = 0
$ws3.Cells.Item(4,22).Value = 0.1
$ws3.Cells.Item(4,23).Value = 150000
$ws3.Cells.Item(4,24).Value = 0
$ws3.Cells.Item(4,25).Value = 0

# Updated "Tổng" row, now row 5
$ws3.Cells.Item(5,2).Value = 3
$ws3.Cells.Item(5,3).Value = 6500000
$ws3.Cells.Item(5,11).Value = 59000000
$ws3.Cells.Item(5,13).Value = 8000000
$ws3.Cells.Item(5,14).Value = 67000000
$ws3.Cells.Item(5,15).Value = 55000000
$ws3.Cells.Item(5,19).Value = 670000
$ws3.Cells.Item(5,23).Value = 650000

# ---------------------------------------------------------------------
# Sheet "Lương"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Lương")

$ws4.Cells.Item(2,2).Value = 8
$ws4.Cells.Item(3,2).Value = 280000
$ws4.Cells.Item(4,2).Value = 952380.9523809524
$ws4.Cells.Item(11,2).Value = 1320000
$ws4.Cells.Item(15,2).Value = 952380.9523809524
$ws4.Cells.Item(26,2).Value = 952380.9523809524
$ws4.Cells.Item(29,2).Value = 1800000
$ws4.Cells.Item(35,2).Value = 3452380.952380952
$ws4.Cells.Item(36,2).Value = 952380.9523809524
$ws4.Cells.Item(37,2).Value = 2752380.952380952
$ws4.Cells.Item(38,2).Value = 7157142.857142857
